$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Trimestre (column C) for all data rows 2-10 to the new quarter.
# Force text format first so Excel keeps this as a literal string instead
# of auto-converting it to a date serial value, then restore the original
# (unstyled / Normal) cell style so formatting is unchanged.
$ws.Range("C2:C10").NumberFormat = "@"
$ws.Range("C2:C10").Value = "01/04/2024"
$ws.Range("C2:C10").Style = "Normal"

# Row 2: Rondônia
$ws.Range("A2").Value = "Rondônia"
$ws.Range("D2").Value = 96.81818181818181
$ws.Range("E2").Value = "1º"

# Row 3: Santa Catarina
$ws.Range("A3").Value = "Santa Catarina"
$ws.Range("D3").Value = 96.80926494918459
$ws.Range("E3").Value = "2º"

# Row 4: Mato Grosso
$ws.Range("A4").Value = "Mato Grosso"
$ws.Range("D4").Value = 96.76767676767678
$ws.Range("E4").Value = "3º"

# Row 5: Mato Grosso do Sul
$ws.Range("A5").Value = "Mato Grosso do Sul"
$ws.Range("D5").Value = 96.18473895582329
$ws.Range("E5").Value = "4º"

# Row 6: Tocantins
$ws.Range("A6").Value = "Tocantins"
$ws.Range("D6").Value = 95.6949569495695
$ws.Range("E6").Value = "5º"

# Row 7: Paraná
$ws.Range("A7").Value = "Paraná"
$ws.Range("D7").Value = 95.56085918854416
$ws.Range("E7").Value = "6º"

# Row 8: Sergipe
$ws.Range("D8").Value = 90.95792300805729
$ws.Range("E8").Value = "21º"

# Row 9: Nordeste
$ws.Range("D9").Value = 90.57060153059201

# Row 10: Brasil
$ws.Range("D10").Value = 93.10426800277951
